$wb = $excel.ActiveWorkbook

# Insert the new "SCMCreds" worksheet right after "UserPageData" (i.e. before "Table"),
# matching the sheet order: LoginPage, LoginPageValidCredentials, UserPageData,
# SCMCreds, Table, UserManagementDropdowns, LoginPageDataProvider, ResetPage.
$afterSheet = $wb.Worksheets.Item("UserPageData")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "SCMCreds"

# Populate the new sheet with the new test-case data.
$newSheet.Range("A1").Value = "Mr"
$newSheet.Range("A2").Value = 1234567890
$newSheet.Range("A3").Value = "PO 45445, NY, 27756"
$newSheet.Range("A4").Value = 34

$newSheet.Columns.Item(1).ColumnWidth = 18.109375

$newSheet.Range("A5").Select() | Out-Null
